# Apply the updated crypto price/volume snapshot (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new display text, and whether the text is a
# number-looking string (e.g. "246.92") that must be forced to remain plain
# text -- otherwise Excel silently re-parses it as a numeric value and the
# exact display string (and any leading/trailing padding) would be lost.
$updates = @(
    @{ Cell = "D2"; Text = '26.494.58'; ForceText = $false }
    @{ Cell = "E2"; Text = '  -0.20%  '; ForceText = $false }
    @{ Cell = "D3"; Text = '1.734.01'; ForceText = $false }
    @{ Cell = "E3"; Text = '  -0.37%  '; ForceText = $false }
    @{ Cell = "E4"; Text = '  +0.10%  '; ForceText = $false }
    @{ Cell = "D5"; Text = '246.92'; ForceText = $true }
    @{ Cell = "E5"; Text = '  +0.55%  '; ForceText = $false }
    @{ Cell = "D7"; Text = '0.4892'; ForceText = $true }
    @{ Cell = "E7"; Text = '  +1.71%  '; ForceText = $false }
    @{ Cell = "E8"; Text = '  -0.66%  '; ForceText = $false }
    @{ Cell = "D9"; Text = '0.06221'; ForceText = $true }
    @{ Cell = "E9"; Text = '  -0.33%  '; ForceText = $false }
    @{ Cell = "D10"; Text = '1.729.60'; ForceText = $false }
    @{ Cell = "E10"; Text = '  -0.62%  '; ForceText = $false }
    @{ Cell = "D11"; Text = '0.07029'; ForceText = $true }
    @{ Cell = "E11"; Text = '  -1.44%  '; ForceText = $false }
    @{ Cell = "E12"; Text = '  -1.26%  '; ForceText = $false }
    @{ Cell = "D13"; Text = '4.588'; ForceText = $true }
    @{ Cell = "E13"; Text = '  +1.05%  '; ForceText = $false }
    @{ Cell = "D14"; Text = '0.6082'; ForceText = $true }
    @{ Cell = "E14"; Text = '  -2.20%  '; ForceText = $false }
    @{ Cell = "D15"; Text = '77.35'; ForceText = $true }
    @{ Cell = "E15"; Text = '  +0.26%  '; ForceText = $false }
    @{ Cell = "E16"; Text = '  +0.05%  '; ForceText = $false }
    @{ Cell = "D17"; Text = '0.000007414'; ForceText = $true }
    @{ Cell = "E17"; Text = '  +7.48%  '; ForceText = $false }
    @{ Cell = "D18"; Text = '26.496.01'; ForceText = $false }
    @{ Cell = "E18"; Text = '  -0.19%  '; ForceText = $false }
    @{ Cell = "E19"; Text = '  +0.05%  '; ForceText = $false }
    @{ Cell = "E20"; Text = '  -2.17%  '; ForceText = $false }
    @{ Cell = "D21"; Text = '1.952.43'; ForceText = $false }
    @{ Cell = "E21"; Text = '  -0.46%  '; ForceText = $false }
    @{ Cell = "E22"; Text = '  -0.64%  '; ForceText = $false }
    @{ Cell = "D23"; Text = '8.739'; ForceText = $true }
    @{ Cell = "E23"; Text = '  -2.07%  '; ForceText = $false }
    @{ Cell = "D24"; Text = '5.228'; ForceText = $true }
    @{ Cell = "E24"; Text = '  -2.26%  '; ForceText = $false }
    @{ Cell = "D25"; Text = '140.94'; ForceText = $true }
    @{ Cell = "E25"; Text = '  +4.00%  '; ForceText = $false }
    @{ Cell = "D26"; Text = '15.42'; ForceText = $true }
    @{ Cell = "E26"; Text = '  +0.20%  '; ForceText = $false }
    @{ Cell = "D27"; Text = '1.416'; ForceText = $true }
    @{ Cell = "E27"; Text = '  -0.33%  '; ForceText = $false }
    @{ Cell = "D28"; Text = '1.770'; ForceText = $true }
    @{ Cell = "E28"; Text = '  -2.42%  '; ForceText = $false }
    @{ Cell = "D29"; Text = '107.77'; ForceText = $true }
    @{ Cell = "E29"; Text = '  +0.81%  '; ForceText = $false }
    @{ Cell = "D30"; Text = '4.018'; ForceText = $true }
    @{ Cell = "E30"; Text = '  +0.33%  '; ForceText = $false }
    @{ Cell = "D31"; Text = '0.08027'; ForceText = $true }
    @{ Cell = "E31"; Text = '  +1.62%  '; ForceText = $false }
    @{ Cell = "D32"; Text = '3.693'; ForceText = $true }
    @{ Cell = "E32"; Text = '  -1.39%  '; ForceText = $false }
    @{ Cell = "D33"; Text = '0.04557'; ForceText = $true }
    @{ Cell = "E33"; Text = '  -0.72%  '; ForceText = $false }
    @{ Cell = "E34"; Text = '  +0.08%  '; ForceText = $false }
    @{ Cell = "D36"; Text = '1.005'; ForceText = $true }
    @{ Cell = "E36"; Text = '  +0.39%  '; ForceText = $false }
    @{ Cell = "D37"; Text = '0.6333'; ForceText = $true }
    @{ Cell = "E37"; Text = '  -0.67%  '; ForceText = $false }
    @{ Cell = "D38"; Text = '0.8951'; ForceText = $true }
    @{ Cell = "D39"; Text = '2.017'; ForceText = $true }
    @{ Cell = "E39"; Text = '  +1.04%  '; ForceText = $false }
    @{ Cell = "D40"; Text = '2.397'; ForceText = $true }
    @{ Cell = "E40"; Text = '  -1.48%  '; ForceText = $false }
    @{ Cell = "D41"; Text = '1.004'; ForceText = $true }
    @{ Cell = "E41"; Text = '  -0.08%  '; ForceText = $false }
    @{ Cell = "D42"; Text = '0.01499'; ForceText = $true }
    @{ Cell = "E42"; Text = '  -1.15%  '; ForceText = $false }
    @{ Cell = "D43"; Text = '101.71'; ForceText = $true }
    @{ Cell = "E43"; Text = '  -9.22%  '; ForceText = $false }
    @{ Cell = "D44"; Text = '5.397'; ForceText = $true }
    @{ Cell = "E44"; Text = '  -6.04%  '; ForceText = $false }
    @{ Cell = "D45"; Text = '0.3886'; ForceText = $true }
    @{ Cell = "E45"; Text = '  -0.93%  '; ForceText = $false }
    @{ Cell = "D46"; Text = '6.916'; ForceText = $true }
    @{ Cell = "E46"; Text = '  -0.68%  '; ForceText = $false }
    @{ Cell = "E47"; Text = '  -1.31%  '; ForceText = $false }
    @{ Cell = "D48"; Text = '0.05397'; ForceText = $true }
    @{ Cell = "E48"; Text = '  +1.20%  '; ForceText = $false }
    @{ Cell = "D49"; Text = '7.807'; ForceText = $true }
    @{ Cell = "E49"; Text = '  -0.99%  '; ForceText = $false }
    @{ Cell = "D50"; Text = '30.45'; ForceText = $true }
    @{ Cell = "E50"; Text = '  -1.38%  '; ForceText = $false }
    @{ Cell = "E51"; Text = '  -0.19%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe is Excel's standard text-entry marker: it forces the
        # value to be stored as a literal string instead of being parsed as a
        # number, without changing the apostrophe itself being part of the value.
        $range.Value = "'" + $u.Text
        # The text-entry marker leaves a "quote prefix" flag on the cell style;
        # reset the style so the cell keeps its original (default) formatting.
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Text
    }
}
